$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Add new headers for columns I and J
$ws.Range("I1").Value = "I0"
$ws.Range("J1").Value = "IF"

# Copy the header formatting (bold font, thin border, centered/top aligned)
# from the existing H1 header cell onto the two new header cells.
$ws.Range("H1").Copy()
$ws.Range("I1:J1").PasteSpecial(-4122)

# Values for columns I (I0) and J (IF), rows 2 through 45
$values = @(
    @(5, 5),
    @(7, 8),
    @(3, 4),
    @(6, 7),
    @(7, 7),
    @(7, 8),
    @(7, 7),
    @(7, 9),
    @(4, 5),
    @(5, 7),
    @(2, 4),
    @(9, 9),
    @(1, 3),
    @(5, 7),
    @(5, 7),
    @(4, 5),
    @(1, 3),
    @(10, 10),
    @(5, 6),
    @(5, 5),
    @(4, 5),
    @(7, 7),
    @(1, 1),
    @(6, 6),
    @(7, 8),
    @(6, 6),
    @(5, 7),
    @(5, 7),
    @(5, 7),
    @(7, 7),
    @(6, 6),
    @(6, 7),
    @(8, 8),
    @(8, 8),
    @(5, 7),
    @(7, 7),
    @(4, 6),
    @(5, 7),
    @(8, 8),
    @(6, 7),
    @(8, 8),
    @(3, 3),
    @(1, 2),
    @(1, 2)
)

for ($idx = 0; $idx -lt $values.Length; $idx++) {
    $row = $idx + 2
    $pair = $values[$idx]
    $ws.Cells.Item($row, 9).Value = $pair[0]
    $ws.Cells.Item($row, 10).Value = $pair[1]
}
